# Apply updated crypto price/volume figures to Sheet1 (rows 2-51)
# Values are stored as plain text in the source data (inline strings),
# so force text number format before assigning to avoid Excel auto-converting
# number-looking strings (e.g. "1.002", "314.52") into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.610.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.52"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4241"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3640"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.12"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07294"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8906"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.70"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.811.97"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.573"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.339"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06893"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.30"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008884"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.611.93"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.988"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.59"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.059.28"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.924"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.20"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "122.66"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.283"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.894"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08936"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7703"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.579"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.919"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.096"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9999"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.101"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05377"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01941"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.826"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.912"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5101"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.280"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06596"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4736"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.38"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.05"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.637"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.64%  "
